$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = "ModifiedFollowing"
$ws.Range("AF2").Value = "ModifiedFollowing"
$ws.Range("U2").Value = "Act/365F"
$ws.Range("AH2").Value = "Act/365F"

$ws.Columns.AutoFit()

$ws.Range("U2,AH2").Select()
